$d = $word.ActiveDocument

# The certificate template's second drawing ("Text Box 2") currently has
# its outline switched off (<a:ln><a:noFill/></a:ln>). Turn the outline
# back on so the box is drawn with a visible 0.75pt border (matching the
# document's "Text 1" theme color).
$shape = $d.Shapes.Item("Text Box 2")
$shape.Line.Visible = $true
